# Correcting Relevance Markers Appenzeller-Herzog (2019) - van Dis (2020)
# Updates metrics values on row 3 (metrics_sim_with_priors.json) of the
# active worksheet to reflect the corrected simulation results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value = 0.9012388162422574
$ws.Range("I3").Value = 0.004344516063543623
$ws.Range("K3").Value = 34.1764705882353

$ws.Range("Q3").Value = 4
$ws.Range("R3").Value = 5
$ws.Range("S3").Value = 12
$ws.Range("T3").Value = 26
$ws.Range("U3").Value = 45
$ws.Range("V3").Value = 5791
$ws.Range("W3").Value = 5790
$ws.Range("X3").Value = 5783
$ws.Range("Y3").Value = 5769
$ws.Range("Z3").Value = 5750

$ws.Range("AF3").Value = 0.99931
$ws.Range("AG3").Value = 0.9991370000000001
$ws.Range("AH3").Value = 0.997929
$ws.Range("AI3").Value = 0.995513
$ws.Range("AJ3").Value = 0.992235
